$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under the duplicate_image_filename column (column E)
# for the existing data rows 2 through 21.
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
